$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 51, pushing existing rows 51..173 down to 52..174.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new data record.
$ws.Range("A51").Value = 8
$ws.Range("B51").Value = "Terminal La Palmera de La Serena"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44519
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112003
$ws.Range("G51").Value = "Ajo"
$ws.Range("H51").Value = "Chino"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 600
$ws.Range("K51").Value = 19000
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = 19500
$ws.Range("N51").Value = "`$/caja 10 kilos"
$ws.Range("O51").Value = "China"
$ws.Range("P51").Value = 1950
$ws.Range("Q51").Value = 10
$ws.Range("R51").Value = "Hortaliza"

# Match the date-number-format style used by the other rows' date column.
$ws.Range("D51").NumberFormat = $ws.Range("D52").NumberFormat
